$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "65.176.15"
$ws.Range("E2").Value = "  -0.69%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.534.18"
$ws.Range("E3").Value = "  +2.63%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# Row 5 - BNB
$ws.Range("D5").Value = "600.44"
$ws.Range("E5").Value = "  +1.39%  "

# Row 6 - Solana
$ws.Range("D6").Value = "138.81"
$ws.Range("E6").Value = "  +0.84%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.532.89"
$ws.Range("E7").Value = "  +2.63%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.09%  "

# Row 9 - XRP
$ws.Range("D9").Value = "0.488"
$ws.Range("E9").Value = "  -2.66%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.80%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "6.90"
$ws.Range("E11").Value = "  -5.92%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "0.390"
$ws.Range("E12").Value = "  +2.72%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "4.138.99"
$ws.Range("E13").Value = "  +2.83%  "

# Row 14 - ShibaInu
$ws.Range("D14").Value = "0.0000185"
$ws.Range("E14").Value = "  +1.80%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "27.11"
$ws.Range("E15").Value = "  +2.21%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.540.60"
$ws.Range("E16").Value = "  +3.48%  "

# Row 17 - TRON
$ws.Range("E17").Value = "  +1.54%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "65.269.93"
$ws.Range("E18").Value = "  -0.46%  "

# Row 20 - Polkadot
$ws.Range("D20").Value = "5.92"
$ws.Range("E20").Value = "  +0.47%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "14.31"
$ws.Range("E21").Value = "  +3.97%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "393.30"
$ws.Range("E22").Value = "  -0.23%  "

# Row 23 - Polygon
$ws.Range("D23").Value = "0.572"
$ws.Range("E23").Value = "  +3.21%  "

# Row 24 - WrappedeETH
$ws.Range("D24").Value = "3.675.10"
$ws.Range("E24").Value = "  +2.51%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "73.80"
$ws.Range("E25").Value = "  +0.58%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.12%  "

# Row 27 - PEPE
$ws.Range("D27").Value = "0.0000115"
$ws.Range("E27").Value = "  +7.45%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "7.72"
$ws.Range("E28").Value = "  +7.46%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  -0.16%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +1.77%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "8.20"
$ws.Range("E31").Value = "  -1.03%  "

# Row 32 - RenzoRestakedETH
$ws.Range("D32").Value = "3.545.96"
$ws.Range("E32").Value = "  +2.80%  "

# Row 33 - USDe
$ws.Range("E33").Value = "  -0.01%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "23.78"
$ws.Range("E34").Value = "  +3.32%  "

# Row 35 - Kaspa
$ws.Range("E35").Value = "  -0.93%  "

# Row 36 - Fetch.AI
$ws.Range("D36").Value = "1.27"
$ws.Range("E36").Value = "  +7.71%  "

# Row 37 - Aptos
$ws.Range("E37").Value = "  -0.04%  "

# Row 38 - Monero
$ws.Range("D38").Value = "168.72"
$ws.Range("E38").Value = "  -2.29%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  +3.67%  "

# Row 40 - NEARProtocol
$ws.Range("D40").Value = "4.94"
$ws.Range("E40").Value = "  +2.57%  "

# Row 41 - Hedera
$ws.Range("D41").Value = "0.0805"
$ws.Range("E41").Value = "  +4.62%  "

# Row 42 - Mantle
$ws.Range("D42").Value = "0.824"
$ws.Range("E42").Value = "  -0.26%  "

# Row 43 - EnergySwap
$ws.Range("D43").Value = "26.31"
$ws.Range("E43").Value = "  +14.00%  "

# Row 44 - OKB
$ws.Range("D44").Value = "42.80"
$ws.Range("E44").Value = "  -2.26%  "

# Row 45 - FirstDigitalUSD
$ws.Range("E45").Value = "  +0.11%  "

# Row 46 - Filecoin
$ws.Range("E46").Value = "  -0.19%  "

# Row 47 - was Stacks, now ONDO
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "1.19"
$ws.Range("E47").Value = "  +6.73%  "

# Row 48 - was ONDO, now Stacks
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "1.67"
$ws.Range("E48").Value = "  +2.44%  "

# Row 49 - Maker
$ws.Range("D49").Value = "2.418.15"
$ws.Range("E49").Value = "  +9.14%  "

# Row 50 - Cosmos
$ws.Range("E50").Value = "  +2.96%  "

# Row 51 - was LidoDAOToken, now Bittensor
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "299.10"
$ws.Range("E51").Value = "  +6.20%  "
